# Daily attendance processing - 2025-12-06 23:46:57
# Reorders the comma-separated "Recorded By" values in column G so that the
# automated/backup identities are listed before the plain "System" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact value replacements observed for column G ("Recorded By").
$map = @{
    "dnasr281@gmail.com, System"           = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"          = "backup@backdoor.com, System"
    "system, System, backup@backdoor.com"  = "backup@backdoor.com, System, system"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
